# The bioSample worksheet tracked an "inductionDelay" boolean-ish flag in
# column H as the number 0 for every data row (rows 2-37). This pass
# switches that flag over to the literal text "None" instead, matching the
# rest of the "major accuracy cleaning" pass on this sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 37; $r++) {
    $ws.Cells.Item($r, 8).Value = "None"
}

# Leave the selection where the editor's cursor ended up after making the
# column-H edits.
$ws.Range("H2:H37").Select()
